# Fruta / hortaliza, semanal
# Insert a new weekly record at row 51, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51 (existing rows 51..124 shift to 52..125)
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly record
$ws.Cells.Item(51, 1).Value = 5
$ws.Cells.Item(51, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(51, 3).Value = "Maule"
$ws.Cells.Item(51, 4).Value = 44483
$ws.Cells.Item(51, 5).Value = 7
$ws.Cells.Item(51, 6).Value = 100112017
$ws.Cells.Item(51, 7).Value = "Apio"
$ws.Cells.Item(51, 8).Value = "Americana (o)"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 600
$ws.Cells.Item(51, 11).Value = 7000
$ws.Cells.Item(51, 12).Value = 7000
$ws.Cells.Item(51, 13).Value = 7000
$ws.Cells.Item(51, 14).Value = "`$/docena de matas"
$ws.Cells.Item(51, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(51, 16).Value = 1167
$ws.Cells.Item(51, 17).Value = 6
$ws.Cells.Item(51, 18).Value = "Hortaliza"
